# Apply cryptos list update (price/volume refresh) per commit:
# "Updated cryptos list on Wed Aug  2 03:26:45 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B / C: coin name + link swaps/updates (plain text, safe to assign directly) ---
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

# --- Column E: volume % strings (already padded w/ spaces, Excel keeps them textual) ---
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("E13").Value = "  +2.99%  "
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("E15").Value = "  +3.51%  "
$ws.Range("E16").Value = "  +3.32%  "
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +3.59%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("E26").Value = "  +2.93%  "
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("E35").Value = "  +6.18%  "
$ws.Range("E36").Value = "  +1.97%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  +10.14%  "
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("E41").Value = "  +6.04%  "
$ws.Range("E42").Value = "  +13.04%  "
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +3.69%  "
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("E48").Value = "  +3.95%  "
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("E50").Value = "  +4.09%  "
$ws.Range("E51").Value = "  +2.26%  "

# --- Column D: price strings. Many look numeric ("0.7015", "5.185", ...) and would be
# auto-converted to real numbers by a plain .Value assignment, which the source file
# never used (every price is stored as literal text). Stage the full D2:D51 text in an
# off-sheet helper column formatted as Text, bulk copy/paste-special (values only) over
# D2:D51 so cell text survives untouched, then wipe the helper column completely.
$helper = $ws.Range("ZZ2:ZZ51")
$helper.NumberFormat = "@"
$ws.Range("ZZ2").Value = "29.766.08"
$ws.Range("ZZ3").Value = "1.865.65"
$ws.Range("ZZ4").Value = "1.0000"
$ws.Range("ZZ5").Value = "246.58"
$ws.Range("ZZ6").Value = "0.7015"
$ws.Range("ZZ7").Value = "1.000"
$ws.Range("ZZ8").Value = "0.07780"
$ws.Range("ZZ9").Value = "0.3085"
$ws.Range("ZZ10").Value = "23.81"
$ws.Range("ZZ11").Value = "0.07845"
$ws.Range("ZZ12").Value = "5.185"
$ws.Range("ZZ13").Value = "92.94"
$ws.Range("ZZ14").Value = "1.854.29"
$ws.Range("ZZ15").Value = "0.6963"
$ws.Range("ZZ16").Value = "6.653"
$ws.Range("ZZ17").Value = "29.755.89"
$ws.Range("ZZ18").Value = "0.000008399"
$ws.Range("ZZ19").Value = "2.116.36"
$ws.Range("ZZ20").Value = "244.19"
$ws.Range("ZZ21").Value = "12.84"
$ws.Range("ZZ22").Value = "1.000"
$ws.Range("ZZ23").Value = "7.664"
$ws.Range("ZZ24").Value = "1.001"
$ws.Range("ZZ25").Value = "0.1519"
$ws.Range("ZZ26").Value = "8.977"
$ws.Range("ZZ27").Value = "160.50"
$ws.Range("ZZ28").Value = "18.43"
$ws.Range("ZZ29").Value = "1.550"
$ws.Range("ZZ30").Value = "4.286"
$ws.Range("ZZ31").Value = "4.212"
$ws.Range("ZZ32").Value = "1.201"
$ws.Range("ZZ33").Value = "0.05113"
$ws.Range("ZZ34").Value = "0.7908"
$ws.Range("ZZ35").Value = "1.926"
$ws.Range("ZZ36").Value = "1.163"
$ws.Range("ZZ37").Value = "2.703"
$ws.Range("ZZ38").Value = "1.342.93"
$ws.Range("ZZ39").Value = "0.01885"
$ws.Range("ZZ40").Value = "2.755"
$ws.Range("ZZ41").Value = "0.9661"
$ws.Range("ZZ42").Value = "6.049"
$ws.Range("ZZ43").Value = "106.86"
$ws.Range("ZZ44").Value = "1.000"
$ws.Range("ZZ45").Value = "9.813"
$ws.Range("ZZ46").Value = "2.016.00"
$ws.Range("ZZ47").Value = "0.00000000125"
$ws.Range("ZZ48").Value = "65.52"
$ws.Range("ZZ49").Value = "0.5196"
$ws.Range("ZZ50").Value = "1.793"
$ws.Range("ZZ51").Value = "7.043"
$helper.Copy()
$ws.Range("D2:D51").PasteSpecial(-4163)
$helper.Clear()

